$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-25 09:48:21"
$wsZh.Range("H4").Value = "2016-03-25 09:49:07"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-25 09:48:32"
$wsDe.Range("H4").Value = "2016-03-25 09:49:23"
